$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Split the single paragraph into two, right before "האם main_general_client"
#    so that paragraph 1 ends with "שאלות למרי – " and paragraph 2 carries the
#    rest of the original content (including the trailing bookmarks).
# ---------------------------------------------------------------------------
$splitRng = $d.Content
$splitRng.Find.Execute(" האם main_general_client") | Out-Null
$splitPoint = $d.Range($splitRng.Start, $splitRng.Start)
$splitPoint.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 2) Remove the leading "האם " word (keep the single leading space) from the
#    start of the (new) second paragraph, turning " האם main_general_client"
#    into " main_general_client".
# ---------------------------------------------------------------------------
$rngAham = $d.Content
$rngAham.Find.Execute("האם main_general_client", $false, $false, $false, $false, $false, $true, 1, $false, "main_general_client", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Insert a new, empty paragraph between paragraph 1 and paragraph 2.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 4) Insert the new lead-in text ("הערות לעצמי – כמה דברים: הסתגלות לselfב ")
#    right before "main_general_client" at the start of paragraph 3.
# ---------------------------------------------------------------------------
$rngLead = $d.Content
$rngLead.Find.Execute("main_general_client") | Out-Null
$leadPoint = $d.Range($rngLead.Start, $rngLead.Start)
$leadPoint.InsertBefore("הערות לעצמי – כמה דברים: הסתגלות לself ב")

# ---------------------------------------------------------------------------
# 5) Rename main_general_client -> main_user_client
# ---------------------------------------------------------------------------
$rngRename = $d.Content
$rngRename.Find.Execute("main_general_client", $false, $false, $false, $false, $false, $true, 1, $false, "main_user_client", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Replace the remainder of the old sentence with the new one.
# ---------------------------------------------------------------------------
$rngTail = $d.Content
$rngTail.Find.Execute(" פותח ישר שרת קבצים או האם יש דרך אחרת לדוגמה כמו דרך השרת. בדיקות אורך של vars האם להפיל את התוכנית\לקוח?", $false, $false, $false, $false, $false, $true, 1, $false, " בנוסף גם יש את העניין של שילוב הקבצים שמתקבל במשתמש, נצטרך לשנות את בקשת הקבצים בשרת ובלקוח הכללי לשנות את שליחת הקבצים ובנוסף כאשר המשתמש מקבל את הקבצים יש לשנות גם את דרך הקבלה באמצעות שינוי המילון", 2) | Out-Null

Write-Host "Paragraphs:" $d.Paragraphs.Count
